$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value2 = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/level-of-care"
$meta.Range("B3").Value2 = "8.0.0"
$meta.Range("B8").Value2 = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value2 = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Move the "ele-1/ext-1" constraint text from the Extension row (row 2)
# down to the Extension.extension row (row 4); row 2's constraint cell
# becomes empty.
$elements = $wb.Worksheets.Item("Elements")
$constraintText = $elements.Range("AI2").Value2
$elements.Range("AI4").Value2 = $constraintText
$elements.Range("AI2").Value2 = ""
